$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZhCn.Columns.Item(3).AutoFit() | Out-Null
$wsDeDe.Columns.Item(3).AutoFit() | Out-Null
